# TC01_Trials_Filter_Gender-Male: updated query in gender script ctdc
#
# Adds a new "TabName" column (with a "CasesTab" row label) in front of
# the existing query/dbExcel/WebExcel columns, and refreshes the Neo4j
# "query" and "StatQuery" cypher text with the new ctdc queries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift everything one column to the right to make room for the new
# "TabName" column at A.
$ws.Columns("A").Insert()

# New header/label column.
$ws.Range("A1").Value = "TabName"
$ws.Range("A2").Value = "CasesTab"

# Updated cypher queries (same cells/styles as before the insert, just
# shifted into B/C).
$ws.Range("B2").Value = "MATCH (ct:clinical_trial)<--(a:arm)<--(c:case)`n    WHERE c.gender = ""MALE""`nWITH DISTINCT c, a, ct`nRETURN `n    COALESCE(c.case_id, '') AS ``Case ID``,`n    COALESCE(ct.clinical_trial_designation, '') AS ``Trial Code``,`n    COALESCE(a.arm_id, '') AS ``Arm``,`n    COALESCE(a.arm_drug, '') AS ``Arm Treatment``,`n    COALESCE(c.disease, '') AS ``Diagnosis``,`n    COALESCE(c.gender, '') AS ``Gender``,`n    COALESCE(c.race, '') AS ``Race``,`n    COALESCE(c.ethnicity, '') AS ``Ethnicity``"

$ws.Range("C2").Value = "MATCH (s:specimen)-->(c:case)-->(:arm)-->(ct:clinical_trial)`n    WHERE c.gender = ""MALE""`nOPTIONAL MATCH (f:file)-->(:sequencing_assay)-->(:nucleic_acid)-->(s)`nRETURN `n    COUNT(DISTINCT f) AS number_of_files,`n    COUNT(DISTINCT c.case_id) AS number_of_cases,`n    COUNT(DISTINCT ct.clinical_trial_designation) AS number_of_trials"

# Match the new row height / column layout / selection seen in the
# authored workbook.
$ws.Rows("2").RowHeight = 174
$ws.Columns("A").ColumnWidth = 8

$ws.Range("C4").Select() | Out-Null
